# New files updated on 3/01/2017
# Adds three new "wish list" locator/description rows to the first
# worksheet (testDataSheet) of the test-data workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = "MyWishList1"
$ws.Range("B4").Value = "New WishList1"

$ws.Range("A5").Value = "MyWishList2"
$ws.Range("B5").Value = "New WishList2"

$ws.Range("A6").Value = "MyWishList3"
$ws.Range("B6").Value = "New WishList3"

# Keep the sheet's recorded selection in step with the extra rows, as
# Excel does when a user keeps typing values down a column.
$ws.Range("B12").Select()
